# "Creacion de la base de datos sobre un .txt"
#
# The currency-lookup helper table (Tabla1, A1:C9) loses its rows for
# CNY, JPY, BRL and KRW (rows 4-9 were TO/FROM pairs keyed on USD/COP);
# only the EUR/USD and EUR/COP rows survive, so the table shrinks to
# A1:C3. The bank-source table (Tabla2, F1:I10) is untouched - its cell
# content doesn't move, it just renumbers against the shared-string
# table once the four now-unused currency codes fall out of it.
#
# Also tidies two stray formatted-but-empty cells (I11, I12), paints
# B14 with the same look as A14, and leaves the selection on B14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shrink Tabla1 from A1:C9 down to A1:C3 (header + the 2 remaining rows).
$lo = $ws.ListObjects("Tabla1")
$lo.Resize($ws.Range("A1:C3"))

# Clear the TO/FROM values and the URL formula/result for the rows that
# fell out of the table (B4:C9 held CNY/JPY/BRL/KRW pairings; A4:A9 held
# the CONCATENATE formula). The style on column A is left in place.
$ws.Range("B4:C9").ClearContents()
$ws.Range("A4:A9").ClearContents()

# Drop the leftover formatting on I11/I12 - they carried a style but no
# value and are no longer needed.
[void]$ws.Range("I11:I12").Clear()

# B14 picks up the same formatting as A14.
[void]$ws.Range("A14").Copy()
[void]$ws.Range("B14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection on B14, matching where editing finished.
[void]$ws.Range("B14").Select()
